$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Paragraph 2 ("{{?comentariosMemoria}}") gains an explicit
#    pageBreakBefore="false" on its pPr.
# -----------------------------------------------------------------------
$pComentariosMemoria = $d.Paragraphs.Item(2)
$pComentariosMemoria.Format.PageBreakBefore = $false

# -----------------------------------------------------------------------
# 2) Insert a brand-new paragraph right before the "Comentarios de
#    Evaluacion para ..." paragraph. The new paragraph only contains a
#    manual page-break run (and an empty leading run).
# -----------------------------------------------------------------------
$pComentarios = $d.Paragraphs.Item(3)
$pComentarios.Range.InsertParagraphBefore()
$pageBreakPara = $d.Paragraphs.Item(3)

$fragment = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr>' +
                '<w:pStyle w:val="Normal"/>' +
                '<w:pageBreakBefore w:val="false"/>' +
                '<w:widowControl/>' +
                '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
                '<w:overflowPunct w:val="true"/>' +
                '<w:spacing w:before="0" w:after="140"/>' +
                '<w:jc w:val="left"/>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Ubuntu" w:hAnsi="Ubuntu"/>' +
                  '<w:sz w:val="24"/>' +
                  '<w:szCs w:val="24"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r><w:rPr/></w:r>' +
              '<w:r><w:br w:type="page"/></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$pageBreakPara.Range.InsertXML($fragment)

# Touch-up: the InsertXML round-trip drops an explicit SpaceBefore="0";
# restore it via the object model without disturbing the runs.
$pageBreakPara2 = $d.Paragraphs.Item(3)
$pageBreakPara2.Format.SpaceBefore = 0

# -----------------------------------------------------------------------
# 3) "Responsable: {{respon" + "s" + "able}}" runs collapse into a
#    single run "Responsable: {{responsable}}".
# -----------------------------------------------------------------------
$d.Content.Find.Execute("Responsable: {{responsable}}", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "Responsable: {{responsable}}", 2)
